$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: stamp the new rows with the same formatting used by the
# existing analogous rows (row 2 = "boxed" name row with top border,
# row 10 = unboxed name row with a Box-column left border, rows 4-8 =
# the plain sub-rows of a customer block). Copy-with-destination carries
# over cell styles without disturbing any other parts of the sheet.
$ws.Range("A2:K2").Copy($ws.Range("A16"))
$ws.Range("A10:K10").Copy($ws.Range("A17"))
$ws.Range("A4:C4").Copy($ws.Range("A18"))
$ws.Range("A5:C5").Copy($ws.Range("A19"))
$ws.Range("A6:C6").Copy($ws.Range("A20"))
$ws.Range("A7:C7").Copy($ws.Range("A21"))
$ws.Range("A8:C8").Copy($ws.Range("A22"))

# --- Step 2: fill in the real content for the new customer block.

# Row 16 - Box 5
$ws.Range("A16").Value = "Name"
$ws.Range("B16").Value = "John Doe"
$ws.Range("C16").Value = "Box 5"
$ws.Range("D16").Value = "INCH"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 234
$ws.Range("I16").Formula = "=E16*2.54 *F16*2.54 *G16*2.54 /6000"
$ws.Range("J16").Value = "Yes"
$ws.Range("K16").Formula = '=IF(J16="No", 13*MAX(H16,I16), IF(J16="Yes", 14*MAX(H16,I16), "Invalid Input—Yes or No"))'

# Row 17 - Box 6
$ws.Range("A17").Value = "Email"
$ws.Range("B17").Value = "john@gmail.com"
$ws.Range("C17").Value = "Box 6"
$ws.Range("D17").Value = "INCH"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 234
$ws.Range("I17").Formula = "=E17*2.54 *F17*2.54 *G17*2.54 /6000"
$ws.Range("J17").Value = "No"
$ws.Range("K17").Formula = '=IF(J17="No", 13*MAX(H17,I17), IF(J17="Yes", 14*MAX(H17,I17), "Invalid Input—Yes or No"))'

# Row 18 - Cell Number (kept as text, same as the sibling "Cell Number"
# rows elsewhere on the sheet, e.g. B4/B11 -- format as Text first so the
# all-digit string isn't auto-coerced into a Number)
$ws.Range("A18").Value = "Cell Number"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2983748932"

# Row 19 - Delivery Option
$ws.Range("A19").Value = "Delivery Option"
$ws.Range("B19").Value = "Door to Door"

# Row 20 - Wants Insurance
$ws.Range("A20").Value = "Wants Insurance"
$ws.Range("B20").Value = $true

# Row 21 - Total Cost
$ws.Range("A21").Value = "Total Cost"
$ws.Range("B21").Formula = "=SUM(K16:K17)"

# Row 22 - Notes (custom row height, blank note field; B22 stays an empty
# placeholder cell the way the Copy in step 1 left it, so we don't touch it)
$ws.Range("A22").Value = "Notes"
$ws.Range("A22").RowHeight = 20
